# Add team record (Wins/Losses/Ties) columns to the player data sheet.
# The new columns AD/AE/AF get the same header style as the existing
# header row (copied from AC1, which already carries style index 1 -
# bold, centered, thin-bordered), and every data row (2-50) gets the
# same W/L/T values (75/86/0) since the whole roster shares one team
# record for the season.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
# Clone the formatting of the last existing header cell (AC1) onto the
# three new header cells so they pick up the bold/centered/bordered
# header style without minting a brand-new style entry.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows -----------------------------------------------------------
$wins = 75
$losses = 86
$ties = 0

for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}
